$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5486
$ws.Range("K3").Value = 5633
$ws.Range("K4").Value = 1175
$ws.Range("K5").Value = 404
$ws.Range("K6").Value = 6261
$ws.Range("K7").Value = 18959

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 347
$ws.Range("K3").Value = 383
$ws.Range("K6").Value = 426
$ws.Range("K7").Value = 1264

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 143
$ws.Range("K3").Value = 152
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 423

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 299
$ws.Range("K4").Value = 38
$ws.Range("K5").Value = 18
$ws.Range("K7").Value = 813

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 187
$ws.Range("K7").Value = 639

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 140
$ws.Range("K6").Value = 158
$ws.Range("K7").Value = 434

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 79
$ws.Range("K7").Value = 318

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 138
$ws.Range("K8").Value = 1264
$ws.Range("K9").Value = 79
$ws.Range("K11").Value = 361
$ws.Range("K12").Value = 35
$ws.Range("K14").Value = 99
$ws.Range("K15").Value = 192
$ws.Range("K19").Value = 553
$ws.Range("K20").Value = 437
$ws.Range("K21").Value = 62
$ws.Range("K22").Value = 51
$ws.Range("K23").Value = 195
$ws.Range("K29").Value = 1023
$ws.Range("K31").Value = 203
$ws.Range("K33").Value = 813
$ws.Range("K37").Value = 639
$ws.Range("K41").Value = 131
$ws.Range("K42").Value = 706
$ws.Range("K43").Value = 165
$ws.Range("K47").Value = 131
$ws.Range("K48").Value = 240
$ws.Range("K49").Value = 105
$ws.Range("K51").Value = 237
$ws.Range("K52").Value = 494
$ws.Range("K53").Value = 240
$ws.Range("K54").Value = 365
$ws.Range("K57").Value = 74
$ws.Range("K60").Value = 117
$ws.Range("K63").Value = 55
$ws.Range("K64").Value = 122
$ws.Range("K65").Value = 434
$ws.Range("K66").Value = 62
$ws.Range("K67").Value = 721
$ws.Range("K73").Value = 166
$ws.Range("K74").Value = 18
$ws.Range("K79").Value = 481
$ws.Range("K83").Value = 423
$ws.Range("K84").Value = 145
$ws.Range("K85").Value = 892
$ws.Range("K89").Value = 276
$ws.Range("K90").Value = 172
$ws.Range("K91").Value = 212
$ws.Range("K94").Value = 257
$ws.Range("K96").Value = 203
$ws.Range("K97").Value = 151
$ws.Range("K99").Value = 318
$ws.Range("K101").Value = 18959

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 68
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 258
$ws.Range("K5").Value = 17
$ws.Range("K7").Value = 721

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 105

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 90
$ws.Range("K7").Value = 365

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 290
$ws.Range("K3").Value = 367
$ws.Range("K6").Value = 290
$ws.Range("K7").Value = 1023

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 57
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 175
$ws.Range("K7").Value = 553

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 51
$ws.Range("K3").Value = 40
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 187
$ws.Range("K3").Value = 218
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 265
$ws.Range("K7").Value = 706

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 195

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 63
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 99
$ws.Range("K7").Value = 212

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 153
$ws.Range("K7").Value = 481

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 139
$ws.Range("K7").Value = 437

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 257

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 68
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 192

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K6").Value = 120
$ws.Range("K7").Value = 361

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K2").Value = 27
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 42
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K2").Value = 29
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 276

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 172

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 237

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 44
$ws.Range("K4").Value = 22
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K6").Value = 220
$ws.Range("K7").Value = 892

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K6").Value = 181
$ws.Range("K7").Value = 494

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 18
